# Apply cryptos list update (Sat Apr 29 17:56:09 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with full content swap (coin identity moved rows) ---
# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("B16").Style = "Normal"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C16").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("E16").Style = "Normal"


# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TRON"
$ws.Range("B17").Style = "Normal"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C17").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06794"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("E17").Style = "Normal"


# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "FraxShare"
$ws.Range("B39").Style = "Normal"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C39").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.028"
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E39").Style = "Normal"


# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("B40").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C40").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5957"
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("E40").Style = "Normal"


# --- Rows with Price (D) and Volume(1h) (E) updates ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.393.40"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("E2").Style = "Normal"


# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.12"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E3").Style = "Normal"


# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +1.04%  "
$ws.Range("E4").Style = "Normal"


# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.50"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E5").Style = "Normal"


# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E6").Style = "Normal"


# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("E7").Style = "Normal"


# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4066"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E8").Style = "Normal"


# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08195"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E9").Style = "Normal"


# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.017"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("E10").Style = "Normal"


# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.38"
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("E11").Style = "Normal"


# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.919.47"
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E12").Style = "Normal"


# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.023"
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("E13").Style = "Normal"


# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.200"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("E14").Style = "Normal"


# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.87"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E15").Style = "Normal"


# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E18").Style = "Normal"


# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.65"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("E19").Style = "Normal"


# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.428.60"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E21").Style = "Normal"


# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.618"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("E22").Style = "Normal"


# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.71"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E23").Style = "Normal"


# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.149.91"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E25").Style = "Normal"


# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.601"
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.28%  "
$ws.Range("E26").Style = "Normal"


# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.92"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E27").Style = "Normal"


# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.05"
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("E28").Style = "Normal"


# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.100"
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E29").Style = "Normal"


# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.93"
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E30").Style = "Normal"


# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.018"
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("E31").Style = "Normal"


# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09550"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E32").Style = "Normal"


# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.597"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("E33").Style = "Normal"


# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.557"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E34").Style = "Normal"


# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.362"
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("E35").Style = "Normal"


# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02278"
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E36").Style = "Normal"


# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06111"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E37").Style = "Normal"


# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.74"
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.07%  "
$ws.Range("E41").Style = "Normal"


# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1848"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E42").Style = "Normal"


# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.406"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E43").Style = "Normal"


# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.246"
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("E44").Style = "Normal"


# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07590"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("E45").Style = "Normal"


# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.43"
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("E46").Style = "Normal"


# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5560"
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E47").Style = "Normal"


# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.947"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E48").Style = "Normal"


# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.10"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("E49").Style = "Normal"


# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.26"
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.50%  "
$ws.Range("E51").Style = "Normal"


# --- Rows with only Volume(1h) (E) updates ---
# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("E20").Style = "Normal"


# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("E24").Style = "Normal"


# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E38").Style = "Normal"


# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.04%  "
$ws.Range("E50").Style = "Normal"


